# Update "Elapsed Duration(Hrs)" (column G) values that were recalculated
# at a slightly later timestamp when the report was regenerated.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "R1"; Cell = "G2"; Value = "3926:33:31" },
    @{ Sheet = "R1"; Cell = "G3"; Value = "66:06:09" },
    @{ Sheet = "R2"; Cell = "G2"; Value = "12107:57:12" },
    @{ Sheet = "R2"; Cell = "G3"; Value = "3237:40:41" },
    @{ Sheet = "R2"; Cell = "G4"; Value = "475:52:15" },
    @{ Sheet = "R4"; Cell = "G2"; Value = "2953:47:01" },
    @{ Sheet = "R4"; Cell = "G3"; Value = "180:59:16" },
    @{ Sheet = "R5"; Cell = "G2"; Value = "427:46:00" },
    @{ Sheet = "R6"; Cell = "G2"; Value = "68:18:18" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
